$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted for Jengibre (Vega Modelo de Temuco)
# right before the existing row 306, pushing all subsequent rows (old 306-364)
# down by one (to 307-365).
$ws.Rows(306).Insert()

$ws.Range("A306").Value = 10
$ws.Range("B306").Value = "Vega Modelo de Temuco"
$ws.Range("C306").Value = "La Araucanía"
$ws.Range("D306").Value = 45209
$ws.Range("E306").Value = 9
$ws.Range("F306").Value = 100114007
$ws.Range("G306").Value = "Jengibre"
$ws.Range("H306").Value = "Sin especificar"
$ws.Range("I306").Value = "Primera"
$ws.Range("J306").Value = 100
$ws.Range("K306").Value = 23000
$ws.Range("L306").Value = 23000
$ws.Range("M306").Value = 23000
$ws.Range("N306").Value = "$/caja 13 kilos"
$ws.Range("O306").Value = "Perú"
$ws.Range("P306").Value = 1769
$ws.Range("Q306").Value = 13
$ws.Range("R306").Value = "Hortaliza"
